# Update the "Price" (D) and "Volume(1h)" (E) columns with refreshed crypto
# market data, as produced by the scheduled GitHub Actions scraper run.
#
# Values assigned to column D are prefixed with a literal leading apostrophe
# ('value) so Excel stores them as text (matching the original inlineStr
# cells) instead of auto-converting numeric-looking strings (e.g. "187.21",
# "0.530", "19.20") into numbers, which would silently drop significant
# trailing zeros or reinterpret thousands-separator dots.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''75.131.32'
$ws.Range('D3').Value = '''2.676.96'
$ws.Range('E3').Value = '  +9.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''187.21'
$ws.Range('E5').Value = '  +12.04%  '
$ws.Range('D6').Value = '''586.31'
$ws.Range('E6').Value = '  +3.04%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +3.86%  '
$ws.Range('E9').Value = '  +9.98%  '
$ws.Range('D10').Value = '''2.674.50'
$ws.Range('E10').Value = '  +9.05%  '
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '''0.357'
$ws.Range('D13').Value = '''4.72'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('D14').Value = '''3.168.30'
$ws.Range('E14').Value = '  +9.14%  '
$ws.Range('D15').Value = '''74.861.23'
$ws.Range('E15').Value = '  +6.79%  '
$ws.Range('D16').Value = '''0.0000187'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').Value = '''26.53'
$ws.Range('E17').Value = '  +9.55%  '
$ws.Range('D18').Value = '''2.705.40'
$ws.Range('E18').Value = '  +10.47%  '
$ws.Range('D19').Value = '''9.15'
$ws.Range('E19').Value = '  +27.82%  '
$ws.Range('E20').Value = '  +8.84%  '
$ws.Range('D21').Value = '''371.58'
$ws.Range('E21').Value = '  +8.66%  '
$ws.Range('D22').Value = '''2.27'
$ws.Range('E22').Value = '  +12.06%  '
$ws.Range('D23').Value = '''4.08'
$ws.Range('E23').Value = '  +4.57%  '
$ws.Range('E24').Value = '  +3.59%  '
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = '''69.77'
$ws.Range('E26').Value = '  +4.86%  '
$ws.Range('E27').Value = '  +7.79%  '
$ws.Range('D28').Value = '''9.35'
$ws.Range('E28').Value = '  +9.20%  '
$ws.Range('D29').Value = '''2.795.12'
$ws.Range('E29').Value = '  +8.32%  '
$ws.Range('E30').Value = '  +2.18%  '
$ws.Range('D31').Value = '''0.0₃0947'
$ws.Range('E31').Value = '  +10.19%  '
$ws.Range('E32').Value = '  +13.45%  '
$ws.Range('D33').Value = '''520.69'
$ws.Range('E33').Value = '  +12.10%  '
$ws.Range('E34').Value = '  +3.47%  '
$ws.Range('E35').Value = '  +7.48%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '''163.52'
$ws.Range('E37').Value = '  +2.55%  '
$ws.Range('E38').Value = '  +5.01%  '
$ws.Range('D39').Value = '''19.20'
$ws.Range('E39').Value = '  +5.31%  '
$ws.Range('D40').Value = '''19.34'
$ws.Range('E40').Value = '  +1.18%  '
$ws.Range('E42').Value = '  +12.76%  '
$ws.Range('D43').Value = '''169.78'
$ws.Range('E43').Value = '  +26.23%  '
$ws.Range('E44').Value = '  +8.50%  '
$ws.Range('D45').Value = '''1.68'
$ws.Range('E45').Value = '  +9.37%  '
$ws.Range('E46').Value = '  +8.34%  '
$ws.Range('E47').Value = '  +11.04%  '
$ws.Range('D48').Value = '''39.09'
$ws.Range('E48').Value = '  +2.55%  '
$ws.Range('D49').Value = '''0.0840'
$ws.Range('E49').Value = '  +15.59%  '
$ws.Range('E50').Value = '  +6.93%  '
$ws.Range('D51').Value = '''0.530'
$ws.Range('E51').Value = '  +7.72%  '
